# Generate Report for Handoff
# Adds two new localization-status rows (108f1c17-... and f0d0dcd2-...)
# around the existing a12b8e80-... row, on all three worksheets
# (Overview, zh-cn, de-de), including the matching hyperlinks.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Drop existing hyperlinks on this sheet; they get rebuilt below in the
# final row positions (the a12b8e80 row moves from row 3 to row 4).
$ws1.Hyperlinks.Delete()

# Row 2 (1a4d56c5...) is untouched, rewrite it as-is so the hyperlink
# rebuild below is consistent.
$ws1.Cells.Item(2,1).Value = "1a4d56c5-df28-4653-9095-a908722463a7.md"
$ws1.Cells.Item(2,2).Value = "Handed back: in sync with en-US"
$ws1.Cells.Item(2,3).Value = "Handed back: in sync with en-US"
$ws1.Cells.Item(2,4).Value = "2016-03-22 04:35:57"
$ws1.Cells.Item(2,4).NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Row 3 (NEW): 108f1c17-...
$ws1.Cells.Item(3,1).Value = "108f1c17-dbd8-4c6f-8dbd-c999531a167b.md"
$ws1.Cells.Item(3,2).Value = "Ready for handoff"
$ws1.Cells.Item(3,3).Value = "Ready for handoff"
$ws1.Cells.Item(3,4).Value = "2016-03-22 04:37:04"
$ws1.Cells.Item(3,4).NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Row 4 (was row 3): a12b8e80-...
$ws1.Cells.Item(4,1).Value = "a12b8e80-94f3-42d6-a9fd-3916c14d06fb.md"
$ws1.Cells.Item(4,2).Value = "Ready for handoff"
$ws1.Cells.Item(4,3).Value = "Ready for handoff"
$ws1.Cells.Item(4,4).Value = "2016-03-22 04:35:31"
$ws1.Cells.Item(4,4).NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Row 5 (NEW): f0d0dcd2-...
$ws1.Cells.Item(5,1).Value = "f0d0dcd2-debc-444d-9e6e-df05d808c4cc.md"
$ws1.Cells.Item(5,2).Value = "Ready for handoff"
$ws1.Cells.Item(5,3).Value = "Ready for handoff"
$ws1.Cells.Item(5,4).Value = "2016-03-22 04:37:04"
$ws1.Cells.Item(5,4).NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Rebuild the "File Name" hyperlinks for rows 2-5
$ws1.Hyperlinks.Add($ws1.Cells.Item(2,1), "https://github.com/OpenLocalizationTest/oltest/blob/a5ea5b9f4fe4e94a6ba29ed731e04de06076347c/e2e/1a4d56c5-df28-4653-9095-a908722463a7.md", "", "", "1a4d56c5-df28-4653-9095-a908722463a7.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Cells.Item(3,1), "https://github.com/OpenLocalizationTest/oltest/blob/c1e6e73c6a7b6a0a2b0d5a47e2a0f8a2e1c6d9a0/e2e/108f1c17-dbd8-4c6f-8dbd-c999531a167b.md", "", "", "108f1c17-dbd8-4c6f-8dbd-c999531a167b.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Cells.Item(4,1), "https://github.com/OpenLocalizationTest/oltest/blob/b26cc27fa94afb827c025e122865ea4cb68ad4c0/e2e/a12b8e80-94f3-42d6-a9fd-3916c14d06fb.md", "", "", "a12b8e80-94f3-42d6-a9fd-3916c14d06fb.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Cells.Item(5,1), "https://github.com/OpenLocalizationTest/oltest/blob/d4f8b6e1a2c3d4e5f6a7b8c9d0e1f2a3b4c5d6e7/e2e/f0d0dcd2-debc-444d-9e6e-df05d808c4cc.md", "", "", "f0d0dcd2-debc-444d-9e6e-df05d808c4cc.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Hyperlinks.Delete()

# Row 2 (1a4d56c5...) untouched content, rewritten for hyperlink rebuild
$ws2.Cells.Item(2,1).Value = "1a4d56c5-df28-4653-9095-a908722463a7.md"
$ws2.Cells.Item(2,2).Value = ".md"
$ws2.Cells.Item(2,3).Value = "Handed back: in sync with en-US"
$ws2.Cells.Item(2,4).Value = "1a4d56c5-df28-4653-9095-a908722463a7.c4f42a81a10cffb31811f5bde29222eb706e78fb.zh-cn.xlf"
$ws2.Cells.Item(2,5).Value = "2016-03-22 04:35:54"
$ws2.Cells.Item(2,5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Cells.Item(2,6).Value = "1a4d56c5-df28-4653-9095-a908722463a7.md"
$ws2.Cells.Item(2,7).Value = "1a4d56c5-df28-4653-9095-a908722463a7.c4f42a81a10cffb31811f5bde29222eb706e78fb.zh-cn.xlf"
$ws2.Cells.Item(2,8).Value = "2016-03-22 04:36:15"
$ws2.Cells.Item(2,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Cells.Item(2,10).Value = "Include"

# Row 3 (NEW): 108f1c17-...
$ws2.Cells.Item(3,1).Value = "108f1c17-dbd8-4c6f-8dbd-c999531a167b.md"
$ws2.Cells.Item(3,2).Value = ".md"
$ws2.Cells.Item(3,3).Value = "Ready for handoff"
$ws2.Cells.Item(3,4).Value = "108f1c17-dbd8-4c6f-8dbd-c999531a167b.e3ac152b8cf002e06b7d7a053acb9070e053c3af.zh-cn.xlf"
$ws2.Cells.Item(3,5).Value = "2016-03-22 04:36:59"
$ws2.Cells.Item(3,5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Cells.Item(3,8).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(3,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Cells.Item(3,10).Value = "Include"

# Row 4 (was row 3): a12b8e80-...
$ws2.Cells.Item(4,1).Value = "a12b8e80-94f3-42d6-a9fd-3916c14d06fb.md"
$ws2.Cells.Item(4,2).Value = ".md"
$ws2.Cells.Item(4,3).Value = "Ready for handoff"
$ws2.Cells.Item(4,4).Value = "a12b8e80-94f3-42d6-a9fd-3916c14d06fb.337521e6ee6f381ee982562de8fe0339b92972af.zh-cn.xlf"
$ws2.Cells.Item(4,5).Value = "2016-03-22 04:35:27"
$ws2.Cells.Item(4,5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Cells.Item(4,8).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(4,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Cells.Item(4,10).Value = "Include"

# Row 5 (NEW): f0d0dcd2-...
$ws2.Cells.Item(5,1).Value = "f0d0dcd2-debc-444d-9e6e-df05d808c4cc.md"
$ws2.Cells.Item(5,2).Value = ".md"
$ws2.Cells.Item(5,3).Value = "Ready for handoff"
$ws2.Cells.Item(5,4).Value = "f0d0dcd2-debc-444d-9e6e-df05d808c4cc.5c5a99335a49859ded9f69f2f81e57fcb47805c2.zh-cn.xlf"
$ws2.Cells.Item(5,5).Value = "2016-03-22 04:36:59"
$ws2.Cells.Item(5,5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Cells.Item(5,8).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(5,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Cells.Item(5,10).Value = "Include"

# Rebuild hyperlinks: column A (.md source) and column D (handoff xlf) for rows 2-5,
# plus F/G (.md target / handoff xlf target) for row 2 only (unique to 1a4d56c5 row).
$ws2.Hyperlinks.Add($ws2.Cells.Item(2,1), "https://github.com/OpenLocalizationTest/oltest/blob/a5ea5b9f4fe4e94a6ba29ed731e04de06076347c/e2e/1a4d56c5-df28-4653-9095-a908722463a7.md", "", "", "1a4d56c5-df28-4653-9095-a908722463a7.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(2,4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0416e60c0ae884232158b822f0bd52324cca8928/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/1a4d56c5-df28-4653-9095-a908722463a7.c4f42a81a10cffb31811f5bde29222eb706e78fb.zh-cn.xlf", "", "", "1a4d56c5-df28-4653-9095-a908722463a7.c4f42a81a10cffb31811f5bde29222eb706e78fb.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(2,6), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0d2a0deb7e54e206e9cebf08b5c55ce0c2b19478/e2e/1a4d56c5-df28-4653-9095-a908722463a7.md", "", "", "1a4d56c5-df28-4653-9095-a908722463a7.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(2,7), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3d734650df474c2f5af9dd9e026f6e8da0e5d691/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/1a4d56c5-df28-4653-9095-a908722463a7.c4f42a81a10cffb31811f5bde29222eb706e78fb.zh-cn.xlf", "", "", "1a4d56c5-df28-4653-9095-a908722463a7.c4f42a81a10cffb31811f5bde29222eb706e78fb.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Cells.Item(3,1), "https://github.com/OpenLocalizationTest/oltest/blob/c1e6e73c6a7b6a0a2b0d5a47e2a0f8a2e1c6d9a0/e2e/108f1c17-dbd8-4c6f-8dbd-c999531a167b.md", "", "", "108f1c17-dbd8-4c6f-8dbd-c999531a167b.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(3,4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e3ac152b8cf002e06b7d7a053acb9070e053c3af/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/108f1c17-dbd8-4c6f-8dbd-c999531a167b.e3ac152b8cf002e06b7d7a053acb9070e053c3af.zh-cn.xlf", "", "", "108f1c17-dbd8-4c6f-8dbd-c999531a167b.e3ac152b8cf002e06b7d7a053acb9070e053c3af.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Cells.Item(4,1), "https://github.com/OpenLocalizationTest/oltest/blob/b26cc27fa94afb827c025e122865ea4cb68ad4c0/e2e/a12b8e80-94f3-42d6-a9fd-3916c14d06fb.md", "", "", "a12b8e80-94f3-42d6-a9fd-3916c14d06fb.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(4,4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/248f7c2f0def21e92af9b1d47caae424fc2ec8dd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a12b8e80-94f3-42d6-a9fd-3916c14d06fb.337521e6ee6f381ee982562de8fe0339b92972af.zh-cn.xlf", "", "", "a12b8e80-94f3-42d6-a9fd-3916c14d06fb.337521e6ee6f381ee982562de8fe0339b92972af.zh-cn.xlf") | Out-Null

$ws2.Hyperlinks.Add($ws2.Cells.Item(5,1), "https://github.com/OpenLocalizationTest/oltest/blob/d4f8b6e1a2c3d4e5f6a7b8c9d0e1f2a3b4c5d6e7/e2e/f0d0dcd2-debc-444d-9e6e-df05d808c4cc.md", "", "", "f0d0dcd2-debc-444d-9e6e-df05d808c4cc.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(5,4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5c5a99335a49859ded9f69f2f81e57fcb47805c2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f0d0dcd2-debc-444d-9e6e-df05d808c4cc.5c5a99335a49859ded9f69f2f81e57fcb47805c2.zh-cn.xlf", "", "", "f0d0dcd2-debc-444d-9e6e-df05d808c4cc.5c5a99335a49859ded9f69f2f81e57fcb47805c2.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Hyperlinks.Delete()

# Row 2 (1a4d56c5...) untouched content, rewritten for hyperlink rebuild
$ws3.Cells.Item(2,1).Value = "1a4d56c5-df28-4653-9095-a908722463a7.md"
$ws3.Cells.Item(2,2).Value = ".md"
$ws3.Cells.Item(2,3).Value = "Handed back: in sync with en-US"
$ws3.Cells.Item(2,4).Value = "1a4d56c5-df28-4653-9095-a908722463a7.c4f42a81a10cffb31811f5bde29222eb706e78fb.de-de.xlf"
$ws3.Cells.Item(2,5).Value = "2016-03-22 04:35:57"
$ws3.Cells.Item(2,5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Cells.Item(2,6).Value = "1a4d56c5-df28-4653-9095-a908722463a7.md"
$ws3.Cells.Item(2,7).Value = "1a4d56c5-df28-4653-9095-a908722463a7.c4f42a81a10cffb31811f5bde29222eb706e78fb.de-de.xlf"
$ws3.Cells.Item(2,8).Value = "2016-03-22 04:36:21"
$ws3.Cells.Item(2,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Cells.Item(2,10).Value = "Include"

# Row 3 (NEW): 108f1c17-...
$ws3.Cells.Item(3,1).Value = "108f1c17-dbd8-4c6f-8dbd-c999531a167b.md"
$ws3.Cells.Item(3,2).Value = ".md"
$ws3.Cells.Item(3,3).Value = "Ready for handoff"
$ws3.Cells.Item(3,4).Value = "108f1c17-dbd8-4c6f-8dbd-c999531a167b.e3ac152b8cf002e06b7d7a053acb9070e053c3af.de-de.xlf"
$ws3.Cells.Item(3,5).Value = "2016-03-22 04:37:04"
$ws3.Cells.Item(3,5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Cells.Item(3,8).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(3,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Cells.Item(3,10).Value = "Include"

# Row 4 (was row 3): a12b8e80-...
$ws3.Cells.Item(4,1).Value = "a12b8e80-94f3-42d6-a9fd-3916c14d06fb.md"
$ws3.Cells.Item(4,2).Value = ".md"
$ws3.Cells.Item(4,3).Value = "Ready for handoff"
$ws3.Cells.Item(4,4).Value = "a12b8e80-94f3-42d6-a9fd-3916c14d06fb.337521e6ee6f381ee982562de8fe0339b92972af.de-de.xlf"
$ws3.Cells.Item(4,5).Value = "2016-03-22 04:35:31"
$ws3.Cells.Item(4,5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Cells.Item(4,8).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(4,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Cells.Item(4,10).Value = "Include"

# Row 5 (NEW): f0d0dcd2-...
$ws3.Cells.Item(5,1).Value = "f0d0dcd2-debc-444d-9e6e-df05d808c4cc.md"
$ws3.Cells.Item(5,2).Value = ".md"
$ws3.Cells.Item(5,3).Value = "Ready for handoff"
$ws3.Cells.Item(5,4).Value = "f0d0dcd2-debc-444d-9e6e-df05d808c4cc.5c5a99335a49859ded9f69f2f81e57fcb47805c2.de-de.xlf"
$ws3.Cells.Item(5,5).Value = "2016-03-22 04:37:04"
$ws3.Cells.Item(5,5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Cells.Item(5,8).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(5,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Cells.Item(5,10).Value = "Include"

# Rebuild hyperlinks
$ws3.Hyperlinks.Add($ws3.Cells.Item(2,1), "https://github.com/OpenLocalizationTest/oltest/blob/a5ea5b9f4fe4e94a6ba29ed731e04de06076347c/e2e/1a4d56c5-df28-4653-9095-a908722463a7.md", "", "", "1a4d56c5-df28-4653-9095-a908722463a7.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(2,4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bac884cec7d3893b6d263cfeac77ade5a0c4a93d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/1a4d56c5-df28-4653-9095-a908722463a7.c4f42a81a10cffb31811f5bde29222eb706e78fb.de-de.xlf", "", "", "1a4d56c5-df28-4653-9095-a908722463a7.c4f42a81a10cffb31811f5bde29222eb706e78fb.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(2,6), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/242339fc50268eb238151d60d31951374f71f8ae/e2e/1a4d56c5-df28-4653-9095-a908722463a7.md", "", "", "1a4d56c5-df28-4653-9095-a908722463a7.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(2,7), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2a3ee98cdc605013532dd9e9bfbd8f7f4ffa757b/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/1a4d56c5-df28-4653-9095-a908722463a7.c4f42a81a10cffb31811f5bde29222eb706e78fb.de-de.xlf", "", "", "1a4d56c5-df28-4653-9095-a908722463a7.c4f42a81a10cffb31811f5bde29222eb706e78fb.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Cells.Item(3,1), "https://github.com/OpenLocalizationTest/oltest/blob/c1e6e73c6a7b6a0a2b0d5a47e2a0f8a2e1c6d9a0/e2e/108f1c17-dbd8-4c6f-8dbd-c999531a167b.md", "", "", "108f1c17-dbd8-4c6f-8dbd-c999531a167b.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(3,4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e3ac152b8cf002e06b7d7a053acb9070e053c3af/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/108f1c17-dbd8-4c6f-8dbd-c999531a167b.e3ac152b8cf002e06b7d7a053acb9070e053c3af.de-de.xlf", "", "", "108f1c17-dbd8-4c6f-8dbd-c999531a167b.e3ac152b8cf002e06b7d7a053acb9070e053c3af.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Cells.Item(4,1), "https://github.com/OpenLocalizationTest/oltest/blob/b26cc27fa94afb827c025e122865ea4cb68ad4c0/e2e/a12b8e80-94f3-42d6-a9fd-3916c14d06fb.md", "", "", "a12b8e80-94f3-42d6-a9fd-3916c14d06fb.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(4,4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/955920f5f6bb53b09b95332bb20170110c24a545/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a12b8e80-94f3-42d6-a9fd-3916c14d06fb.337521e6ee6f381ee982562de8fe0339b92972af.de-de.xlf", "", "", "a12b8e80-94f3-42d6-a9fd-3916c14d06fb.337521e6ee6f381ee982562de8fe0339b92972af.de-de.xlf") | Out-Null

$ws3.Hyperlinks.Add($ws3.Cells.Item(5,1), "https://github.com/OpenLocalizationTest/oltest/blob/d4f8b6e1a2c3d4e5f6a7b8c9d0e1f2a3b4c5d6e7/e2e/f0d0dcd2-debc-444d-9e6e-df05d808c4cc.md", "", "", "f0d0dcd2-debc-444d-9e6e-df05d808c4cc.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(5,4), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5c5a99335a49859ded9f69f2f81e57fcb47805c2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f0d0dcd2-debc-444d-9e6e-df05d808c4cc.5c5a99335a49859ded9f69f2f81e57fcb47805c2.de-de.xlf", "", "", "f0d0dcd2-debc-444d-9e6e-df05d808c4cc.5c5a99335a49859ded9f69f2f81e57fcb47805c2.de-de.xlf") | Out-Null

Write-Host "Localization status report regenerated for handoff."
